# Update this month's driving-summary input figures for each driver.
# Only the raw input cells are written; every dependent formula cell
# (C5, F7, G7, C8, C10, F10, G10, C12, F13, G13, C14, ...) is left as a
# formula and refreshed by the recalculation below.

$wb = $excel.ActiveWorkbook

$stefan     = $wb.Worksheets.Item("Stefan")
$christiaan = $wb.Worksheets.Item("Christiaan")
$derrick    = $wb.Worksheets.Item("Derrick")

# --- Stefan: raw inputs ---
$stefan.Range("F6").Value = 1097
$stefan.Range("G6").Value = 856
$stefan.Range("C7").Value = 289

# --- Christiaan: raw inputs ---
$christiaan.Range("F6").Value = 817
$christiaan.Range("G6").Value = 895
$christiaan.Range("C7").Value = 234
$christiaan.Range("C13").Value = 1048.71

# --- Derrick: raw inputs ---
$derrick.Range("F6").Value = 490
$derrick.Range("G6").Value = 890

# Recalculate the whole workbook so every dependent formula picks up
# the new inputs.
$excel.CalculateFullRebuild()

# --- Restore the per-sheet selected cell ---
$stefan.Range("G6").Select() | Out-Null
$christiaan.Range("G6").Select() | Out-Null
$derrick.Range("F6").Select() | Out-Null

# --- Christiaan becomes the active (displayed) sheet ---
$christiaan.Activate()

$wb.Save()
